$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.746.08"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "1.874.60"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").Value = "'325.26"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.4592"
$ws.Range("D8").Value = "'0.3861"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.07866"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").Value = "'0.9944"
$ws.Range("E10").Value = "  +3.37%  "
$ws.Range("D11").Value = "'21.70"
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("D12").Value = "1.896.40"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("D13").Value = "'6.982"
$ws.Range("D14").Value = "'5.700"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "'0.06973"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "'88.39"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "'0.00001005"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").Value = "'16.83"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "28.754.08"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").Value = "'5.277"
$ws.Range("D23").Value = "'11.01"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").Value = "2.116.87"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("D26").Value = "'153.39"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "'5.783"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").Value = "'1.953"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "'118.93"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'0.09318"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").Value = "'0.9189"
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("D33").Value = "'5.300"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").Value = "'3.319"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "'0.05749"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("D38").Value = "'0.02074"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").Value = "'7.715"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "'0.5636"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "'0.1788"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").Value = "'9.883"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "'0.07198"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").Value = "'11.70"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'0.5278"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "'2.144"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "'1.118"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.826"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'113.47"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "'2.411"
$ws.Range("D51").Value = "'1.003"
$ws.Range("E51").Value = "  +0.33%  "
